# Auto-generated edit script applying numeric cell changes to the Leve profit tables
# (columns H..N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5321.826
$ws.Range("I43").Value = 3742.6667
$ws.Range("J43").Value = 5558.7
$ws.Range("K43").Value = 3742.6667
$ws.Range("L43").Value = 5558.7
$ws.Range("M43").Value = -3673.6667
$ws.Range("N43").Value = -5696.7
$ws.Range("H88").Value = 6130
$ws.Range("I88").Value = 6130
$ws.Range("K88").Value = 6130
$ws.Range("M88").Value = -5724
$ws.Range("H91").Value = 6130
$ws.Range("I91").Value = 6130
$ws.Range("K91").Value = 6130
$ws.Range("M91").Value = -4726
$ws.Range("H98").Value = 631.2857
$ws.Range("I98").Value = 649.0769
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 649.0769
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 848.9231
$ws.Range("N98").Value = -3396
$ws.Range("H113").Value = 2441.1765
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -9508
$ws.Range("H116").Value = 3095.3462
$ws.Range("I116").Value = 1300
$ws.Range("J116").Value = 3167.16
$ws.Range("K116").Value = 1300
$ws.Range("L116").Value = 3167.16
$ws.Range("M116").Value = 2142
$ws.Range("N116").Value = -10051.16
$ws.Range("H122").Value = 631.2857
$ws.Range("I122").Value = 649.0769
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 1947.2307
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 502.7692999999999
$ws.Range("N122").Value = -6100
$ws.Range("H139").Value = 40157.445
$ws.Range("J139").Value = 49881.668
$ws.Range("L139").Value = 49881.668
$ws.Range("N139").Value = -60161.668
$ws.Range("H141").Value = 1274.2307
$ws.Range("I141").Value = 1161
$ws.Range("K141").Value = 3483
$ws.Range("M141").Value = 1697
# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11365824
$ws.Range("I74").Value = 15626144
$ws.Range("J74").Value = 4969
$ws.Range("K74").Value = 15626144
$ws.Range("L74").Value = 4969
$ws.Range("M74").Value = -15625270
$ws.Range("N74").Value = -6717
$ws.Range("H77").Value = 11365824
$ws.Range("I77").Value = 15626144
$ws.Range("J77").Value = 4969
$ws.Range("K77").Value = 78130720
$ws.Range("L77").Value = 24845
$ws.Range("M77").Value = -78126352
$ws.Range("N77").Value = -33581
# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1456.4783
$ws.Range("I99").Value = 1366.6666
$ws.Range("J99").Value = 1554.4546
$ws.Range("K99").Value = 1366.6666
$ws.Range("L99").Value = 1554.4546
$ws.Range("M99").Value = 131.3334
$ws.Range("N99").Value = -4550.4546
$ws.Range("H126").Value = 1456.4783
$ws.Range("I126").Value = 1366.6666
$ws.Range("J126").Value = 1554.4546
$ws.Range("K126").Value = 4099.9998
$ws.Range("L126").Value = 4663.3638
$ws.Range("M126").Value = -1629.9998
$ws.Range("N126").Value = -9603.363799999999
$ws.Range("H132").Value = 16668805
$ws.Range("I132").Value = 22728892
$ws.Range("K132").Value = 68186676
$ws.Range("M132").Value = -68184146
$ws.Range("H134").Value = 1329.8055
$ws.Range("I134").Value = 1191.3928
$ws.Range("J134").Value = 1814.25
$ws.Range("K134").Value = 3574.1784
$ws.Range("L134").Value = 5442.75
$ws.Range("M134").Value = -1039.1784
$ws.Range("N134").Value = -10512.75
$ws.Range("H140").Value = 45094.145
$ws.Range("J140").Value = 45094.145
$ws.Range("L140").Value = 45094.145
$ws.Range("N140").Value = -55454.145
# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 180
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -11
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 1119.8823
$ws.Range("I34").Value = 223.33333
$ws.Range("J34").Value = 1608.909
$ws.Range("K34").Value = 669.99999
$ws.Range("L34").Value = 4826.727000000001
$ws.Range("M34").Value = -585.99999
$ws.Range("N34").Value = -4994.727000000001
$ws.Range("H39").Value = 500.55554
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 508.26923
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 1524.80769
$ws.Range("M39").Value = -606
$ws.Range("N39").Value = -2112.80769
$ws.Range("H55").Value = 525
$ws.Range("J55").Value = 606.25
$ws.Range("L55").Value = 1818.75
$ws.Range("N55").Value = -2172.75
$ws.Range("H131").Value = 2649.6667
$ws.Range("J131").Value = 4582.6665
$ws.Range("L131").Value = 13747.9995
$ws.Range("N131").Value = -23827.9995
$ws.Range("H138").Value = 4220
$ws.Range("J138").Value = 9644.333000000001
$ws.Range("L138").Value = 28932.999
$ws.Range("N138").Value = -39212.999
# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59950
$ws.Range("J138").Value = 59950
$ws.Range("L138").Value = 59950
$ws.Range("N138").Value = -70230
# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5421.0415
$ws.Range("I7").Value = 5397.5757
$ws.Range("J7").Value = 5472.6665
$ws.Range("K7").Value = 5397.5757
$ws.Range("L7").Value = 5472.6665
$ws.Range("M7").Value = -5285.5757
$ws.Range("N7").Value = -5696.6665
$ws.Range("H40").Value = 5821.316
$ws.Range("I40").Value = 8757.143
$ws.Range("J40").Value = 4108.75
$ws.Range("K40").Value = 8757.143
$ws.Range("L40").Value = 4108.75
$ws.Range("M40").Value = -8621.143
$ws.Range("N40").Value = -4380.75
$ws.Range("H126").Value = 5421.0415
$ws.Range("I126").Value = 5397.5757
$ws.Range("J126").Value = 5472.6665
$ws.Range("K126").Value = 16192.7271
$ws.Range("L126").Value = 16417.9995
$ws.Range("M126").Value = -13722.7271
$ws.Range("N126").Value = -21357.9995
$ws.Range("H132").Value = 9476.941999999999
$ws.Range("I132").Value = 6436.4
$ws.Range("J132").Value = 11757.35
$ws.Range("K132").Value = 19309.2
$ws.Range("L132").Value = 35272.05
$ws.Range("M132").Value = -16779.2
$ws.Range("N132").Value = -40332.05
$ws.Range("H139").Value = 60707.5
$ws.Range("J139").Value = 60707.5
$ws.Range("L139").Value = 60707.5
$ws.Range("N139").Value = -70987.5
# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10038
$ws.Range("I126").Value = 4968
$ws.Range("K126").Value = 14904
$ws.Range("M126").Value = -12434
$ws.Range("H132").Value = 2333.138
$ws.Range("I132").Value = 1694.0834
$ws.Range("J132").Value = 5400.6
$ws.Range("K132").Value = 5082.2502
$ws.Range("L132").Value = 16201.8
$ws.Range("M132").Value = -2552.2502
$ws.Range("N132").Value = -21261.8
